$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "E2" "-0.73%"
Set-TextValue "D3" "27.01"
Set-TextValue "E3" "3.66%"
Set-TextValue "D4" "5.155"
Set-TextValue "E4" "1.08%"
Set-TextValue "D5" "0.05623"
Set-TextValue "E5" "0.51%"
Set-TextValue "D6" "6.476"
Set-TextValue "E6" "-0.23%"
Set-TextValue "D7" "0.8167"
Set-TextValue "D8" "0.8325"
Set-TextValue "E8" "-1.60%"
Set-TextValue "D9" "0.1327"
Set-TextValue "E9" "-1.13%"
Set-TextValue "D10" "0.06923"
Set-TextValue "E10" "-0.41%"
Set-TextValue "D11" "0.02890"
Set-TextValue "E11" "1.38%"
Set-TextValue "D12" "0.09384"
Set-TextValue "E12" "0.08%"
Set-TextValue "D13" "0.001509"
Set-TextValue "E13" "-0.52%"
Set-TextValue "D14" "0.04268"
Set-TextValue "E14" "-8.96%"
Set-TextValue "D15" "0.0005998"
Set-TextValue "E15" "-93.87%"
Set-TextValue "D16" "0.006111"
Set-TextValue "E16" "-1.33%"
Set-TextValue "D17" "3.607"
Set-TextValue "E17" "1.48%"
Set-TextValue "D18" "3.019"
Set-TextValue "E18" "-0.05%"
Set-TextValue "D19" "2.307"
Set-TextValue "E19" "8.91%"
Set-TextValue "D21" "0.03100"
Set-TextValue "E21" "-3.57%"
Set-TextValue "E22" "-2.18%"
Set-TextValue "D23" "3.742"
Set-TextValue "E23" "0.06%"
Set-TextValue "E24" "-0.08%"
Set-TextValue "E25" "-1.96%"
Set-TextValue "E26" "-2.85%"
Set-TextValue "D27" "0.00009796"
Set-TextValue "E27" "2.07%"
Set-TextValue "E28" "-0.45%"
Set-TextValue "D40" "0.03649"
Set-TextValue "E40" "-0.15%"
Set-TextValue "D41" "0.006038"
Set-TextValue "E41" "-1.32%"
Set-TextValue "E42" "-0.14%"
Set-TextValue "D43" "0.002556"
Set-TextValue "E43" "2.27%"
Set-TextValue "D44" "0.008168"
Set-TextValue "E44" "5.18%"
Set-TextValue "D45" "0.00005306"
Set-TextValue "E45" "-0.15%"
Set-TextValue "E46" "-0.01%"
Set-TextValue "D47" "0.1090"
Set-TextValue "E47" "-35.89%"
Set-TextValue "E48" "28.63%"
Set-TextValue "E49" "-0.01%"
Set-TextValue "E50" "-0.01%"
